# Apply the "Added PDF versions to site" edit:
#  1. Remove the "Word version of this document" bullet from the
#     "Additional resources" list (its target is being replaced by a PDF
#     link elsewhere, so the whole list item goes away).
#  2. Tidy the wording in the inflection-points paragraph under "Pitfalls".

$d = $word.ActiveDocument

# --- 1. Remove the "Word version of this document" list paragraph -------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Word version of this document*") {
        $para.Range.Delete()
        break
    }
}

# --- 2. Fix the "inflection points" sentence -----------------------------
$d.Content.Find.Execute(
    "a lesson about inflection points. There is any particular relevance to statistics.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a lesson about inflection points, not about statistics.", 2
) | Out-Null

Write-Output "done"
